$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "customer maintenance page added"
# Row 4's method name is repurposed from the old "verifyValueTypeSettingMap"
# test to the new customer-navigation test.
$ws.Range("A4").Value = "verifyCustomerNavigationAfterLogin"

# Add the new "customer maintenance" data-entry test results (currently
# failing), six rows' worth, mirroring the existing Method/Status/Date
# layout used by the other rows.
for ($i = 5; $i -le 10; $i++) {
    $ws.Cells.Item($i, 1).Value = "verifyCustomerDataEntry"

    $statusCell = $ws.Cells.Item($i, 2)
    $statusCell.Value = "Failed"
    # Existing "Passed" rows are shaded green (ColorIndex 10); shade the new
    # "Failed" rows red (ColorIndex 3) so failures stand out the same way.
    $statusCell.Interior.ColorIndex = 3

    $ws.Cells.Item($i, 3).Value = "28-12-2024"
}

# Column A grew slightly wider to fit the longer method names.
$ws.Columns.Item(1).ColumnWidth = 32.86
